$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 101.14286
$ws.Range("I6").Value = 34.666668
$ws.Range("K6").Value = 104.000004
$ws.Range("M6").Value = 7.999995999999996
$ws.Range("H9").Value = 109.2
$ws.Range("I9").Value = 86.5
$ws.Range("K9").Value = 86.5
$ws.Range("M9").Value = 82.5
$ws.Range("H40").Value = 6704.154
$ws.Range("I40").Value = 5832
$ws.Range("J40").Value = 8099.6
$ws.Range("K40").Value = 5832
$ws.Range("L40").Value = 8099.6
$ws.Range("M40").Value = -5657
$ws.Range("N40").Value = -8449.6
$ws.Range("H41").Value = 898.3
$ws.Range("I41").Value = 775.8889
$ws.Range("K41").Value = 775.8889
$ws.Range("M41").Value = -335.8889
$ws.Range("H43").Value = 2341.75
$ws.Range("I43").Value = 2155.6667
$ws.Range("K43").Value = 2155.6667
$ws.Range("M43").Value = -2086.6667
$ws.Range("H53").Value = 1057.4286
$ws.Range("I53").Value = 624.1429000000001
$ws.Range("K53").Value = 624.1429000000001
$ws.Range("M53").Value = 12.85709999999995
$ws.Range("H55").Value = 4874
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 4874
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 4874
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -5302
$ws.Range("H101").Value = 413.83334
$ws.Range("J101").Value = 477.5
$ws.Range("L101").Value = 1432.5
$ws.Range("N101").Value = -4676.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6597.75
$ws.Range("I2").Value = 1040.4286
$ws.Range("K2").Value = 1040.4286
$ws.Range("M2").Value = -927.4286
$ws.Range("H74").Value = 3262.0527
$ws.Range("I74").Value = 2959.7693
$ws.Range("K74").Value = 2959.7693
$ws.Range("M74").Value = -2085.7693
$ws.Range("H77").Value = 3262.0527
$ws.Range("I77").Value = 2959.7693
$ws.Range("K77").Value = 14798.8465
$ws.Range("M77").Value = -10430.8465
$ws.Range("H97").Value = 1841.5
$ws.Range("I97").Value = 1841.5
$ws.Range("K97").Value = 1841.5
$ws.Range("M97").Value = -1345.5
$ws.Range("H102").Value = 15631062
$ws.Range("J102").Value = 7583.3335
$ws.Range("L102").Value = 7583.3335
$ws.Range("N102").Value = -10827.3335
$ws.Range("H116").Value = 6597.75
$ws.Range("I116").Value = 1040.4286
$ws.Range("K116").Value = 1040.4286
$ws.Range("M116").Value = 1253.5714
$ws.Range("H122").Value = 2830
$ws.Range("I122").Value = 2830
$ws.Range("K122").Value = 8490
$ws.Range("M122").Value = -6040
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6597.75
$ws.Range("I3").Value = 1040.4286
$ws.Range("K3").Value = 1040.4286
$ws.Range("M3").Value = -926.4286
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H105").Value = 30304368
$ws.Range("I105").Value = 30304368
$ws.Range("K105").Value = 30304368
$ws.Range("M105").Value = -30302621
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 866.41174
$ws.Range("J22").Value = 1075.2
$ws.Range("L22").Value = 1075.2
$ws.Range("N22").Value = -1775.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 128.38461
$ws.Range("J12").Value = 175.33333
$ws.Range("L12").Value = 525.99999
$ws.Range("N12").Value = -871.99999
$ws.Range("H17").Value = 44.285713
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 48.333332
$ws.Range("K17").Value = 60
$ws.Range("L17").Value = 144.999996
$ws.Range("M17").Value = 109
$ws.Range("N17").Value = -482.999996
$ws.Range("H33").Value = 33.25
$ws.Range("I33").Value = 38.666668
$ws.Range("J33").Value = 17
$ws.Range("K33").Value = 232.000008
$ws.Range("L33").Value = 102
$ws.Range("M33").Value = 50.99999199999999
$ws.Range("N33").Value = -668
$ws.Range("H38").Value = 437.22223
$ws.Range("I38").Value = 441.875
$ws.Range("K38").Value = 1325.625
$ws.Range("M38").Value = -978.625
$ws.Range("H58").Value = 5695
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H121").Value = 1291.6666
$ws.Range("I121").Value = 437.5
$ws.Range("K121").Value = 1312.5
$ws.Range("M121").Value = -2.5
$ws.Range("H139").Value = 5473.8
$ws.Range("I139").Value = 5473.8
$ws.Range("K139").Value = 16421.4
$ws.Range("M139").Value = -11281.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 413.9091
$ws.Range("I2").Value = 131.625
$ws.Range("J2").Value = 1166.6666
$ws.Range("K2").Value = 131.625
$ws.Range("L2").Value = 1166.6666
$ws.Range("M2").Value = -18.625
$ws.Range("N2").Value = -1392.6666
$ws.Range("H23").Value = 600
$ws.Range("J23").Value = 600
$ws.Range("L23").Value = 600
$ws.Range("N23").Value = -1046
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1500
$ws.Range("K83").Value = 7500
$ws.Range("M83").Value = -2508
$ws.Range("H102").Value = 3804.6155
$ws.Range("I102").Value = 3370.125
$ws.Range("K102").Value = 3370.125
$ws.Range("M102").Value = -1748.125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1684.3334
$ws.Range("J7").Value = 1548
$ws.Range("L7").Value = 1548
$ws.Range("N7").Value = -1772
$ws.Range("H46").Value = 5476.095
$ws.Range("I46").Value = 5499.857
$ws.Range("J46").Value = 5464.2144
$ws.Range("K46").Value = 5499.857
$ws.Range("L46").Value = 5464.2144
$ws.Range("M46").Value = -5311.857
$ws.Range("N46").Value = -5840.2144
$ws.Range("H87").Value = 20000
$ws.Range("J87").Value = 20000
$ws.Range("L87").Value = 20000
$ws.Range("N87").Value = -22246
$ws.Range("H90").Value = 20000
$ws.Range("J90").Value = 20000
$ws.Range("L90").Value = 60000
$ws.Range("N90").Value = -71232
$ws.Range("H122").Value = 5710.625
$ws.Range("I122").Value = 6126.4287
$ws.Range("K122").Value = 18379.2861
$ws.Range("M122").Value = -15929.2861
$ws.Range("H126").Value = 1684.3334
$ws.Range("J126").Value = 1548
$ws.Range("L126").Value = 4644
$ws.Range("N126").Value = -9584
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 17750
$ws.Range("J56").Value = 17750
$ws.Range("L56").Value = 17750
$ws.Range("N56").Value = -19178
$ws.Range("H96").Value = 1320.2
$ws.Range("I96").Value = 1275.25
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 1275.25
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = 97.75
$ws.Range("N96").Value = -4246
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H126").Value = 5348.125
$ws.Range("I126").Value = 928.3333
$ws.Range("K126").Value = 2784.9999
$ws.Range("M126").Value = -314.9998999999998
$ws.Range("H135").Value = 87500
$ws.Range("J135").Value = 87500
$ws.Range("L135").Value = 87500
$ws.Range("N135").Value = -97640
$ws.Range("H136").Value = 3189.9285
$ws.Range("I136").Value = 2196.9
$ws.Range("J136").Value = 5672.5
$ws.Range("K136").Value = 6590.700000000001
$ws.Range("L136").Value = 17017.5
$ws.Range("M136").Value = -4040.700000000001
$ws.Range("N136").Value = -22117.5
